# 修改 /api/order/{id} get 接口
# - URL changes from /api/order/id/{id} to /api/order/{id}
# - request-parameter cell (E36) is removed
# - response-data cell (F36) gains an eContractStatus field
# - row 36 grows taller to fit the extra line
# - the PUT /api/order/{id} endpoint's eContractStatus param description (E37)
#   is expanded with allowed literal values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API")

# B36: endpoint path for "获取订单详情" (GET order detail)
$ws.Range("B36").Value = "/api/order/{id}"

# E36: the old "String id //订单Id" request-param note is dropped entirely
$ws.Range("E36").Clear()

# F36: response payload now also reports the order's confirmation status
$ws.Range("F36").Value = "data:{id:订单Id,
contractNo: 合同编号,
productType: 产品类型（开放式，封闭式，私教课),
validityTimes: 合约有效次数,
contractStart: 合约开始时间,
contractEnd: 合约结束时间,
courseId: 课程Id,
courseName: 课程名，
courseStore:上课门店，
courseStartDate: 课程开始日期,
eContractStatus: 订单的确认状态(true,false)}"

# row 36 needs extra height for the longer F36 text
$ws.Rows.Item(36).RowHeight = 148.5

# E37: describe the accepted literal encodings for the boolean flag
$ws.Range("E37").Value = "{eContractStatus:boolean (0,1,'true','false')}"

# Restore the frozen-header view, scrolled/selected the way the author left it
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("E38").Select()
